$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.171.59"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").Value = "3.881.05"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "482.42"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.22"
$ws.Range("E6").Value = "  -1.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.739"
$ws.Range("E9").Value = "  +2.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000351"
$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.95"
$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.45"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "4.500.54"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "3.919.99"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.17"
$ws.Range("E16").Value = "  -2.81%  "

$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.90"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").Value = "68.202.67"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.94"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.56"
$ws.Range("E22").Value = "  +7.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.72"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.31"
$ws.Range("E24").Value = "  +17.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.70"
$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.67"
$ws.Range("E26").Value = "  +2.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -5.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.06"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  -3.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "718.62"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.44"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("E33").Value = "  +2.60%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0882"
$ws.Range("E34").Value = "  -0.36%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.74"
$ws.Range("E35").Value = "  +5.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  +8.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.67"
$ws.Range("E37").Value = "  -1.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.400"
$ws.Range("E38").Value = "  +17.55%  "

$ws.Range("E39").Value = "  -3.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0494"
$ws.Range("E41").Value = "  +5.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.98"
$ws.Range("E42").Value = "  +8.14%  "

$ws.Range("E43").Value = "  +3.11%  "

$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.142"
$ws.Range("E45").Value = "  +1.08%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.37"
$ws.Range("E46").Value = "  +7.45%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "0.0₆0355"
$ws.Range("E48").Value = "  +28.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.35"
$ws.Range("E49").Value = "  -1.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.09"
$ws.Range("E50").Value = "  -2.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.19"
$ws.Range("E51").Value = "  -2.68%  "
